# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 3 ("R") updates
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("B3").Value = 368
$wsOFF.Range("C3").Value = 259
$wsOFF.Range("D3").Value = 95
$wsOFF.Range("F3").Value = 6

# Sheet "DEF" - row 3 ("R") updates
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("B3").Value = 410
$wsDEF.Range("C3").Value = 296
$wsDEF.Range("D3").Value = 82
$wsDEF.Range("E3").Value = 42
